$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B8").Value = "Computadora"
$ws.Range("C8").Value = "No se han encontrado resultados para la búsqueda."
$ws.Range("C10").Select()
